# Auto-applied edit script for 上海-漫展信息.xlsx
# Updates "想去人数" (F column) counts across all 4 sheets, and refreshes
# rows 38-40 of sheet 1 (展览) with new listing content (gh-pages re-scrape).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1585
$ws.Range("F6").Value = 769
$ws.Range("F7").Value = 685
$ws.Range("F8").Value = 1281
$ws.Range("F9").Value = 2542
$ws.Range("F10").Value = 1329
$ws.Range("F11").Value = 296
$ws.Range("F12").Value = 2297
$ws.Range("F13").Value = 2000
$ws.Range("F15").Value = 6203
$ws.Range("F16").Value = 111
$ws.Range("F17").Value = 1205
$ws.Range("F18").Value = 127
$ws.Range("F19").Value = 1418
$ws.Range("F20").Value = 1314
$ws.Range("F21").Value = 1182
$ws.Range("F22").Value = 97
$ws.Range("F23").Value = 2088
$ws.Range("F25").Value = 679
$ws.Range("F26").Value = 215
$ws.Range("F27").Value = 5223
$ws.Range("F29").Value = 1238
$ws.Range("F31").Value = 3657
$ws.Range("F32").Value = 633
$ws.Range("F33").Value = 1659
$ws.Range("F35").Value = 143
$ws.Range("F36").Value = 265
$ws.Range("F37").Value = 956
$ws.Range("C38").Value = "上海·寻漫岛动漫嘉年华"
$ws.Range("D38").Value = "中山北路3300号4楼L4001号 环球港上海世嘉都市乐园"
$ws.Range("E38").Value = "2024.08.17 10:00-08.18 17:00"
$ws.Range("F38").Value = 38
$ws.Range("G38").Value = 60
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=87628"
$ws.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202407/ePcJqKzI1721099263380.jpeg"
$ws.Range("C39").Value = "上海·火影忍者ONLY·霓虹夏日"
$ws.Range("D39").Value = "沪太路3651弄红光体育运动中心199号 堂颂羽毛球馆"
$ws.Range("E39").Value = "2024.08.17 10:00-08.17 16:00"
$ws.Range("F39").Value = 372
$ws.Range("G39").Value = 89
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=87882"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202406/G4LTUIF51719209197774.jpeg"
$ws.Range("C40").Value = "上海·第六届燃梦BACG PRO动漫嘉年华-我们在燃梦相遇吧！"
$ws.Range("D40").Value = "盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)"
$ws.Range("E40").Value = "2024.08.17 11:00-08.18 16:00"
$ws.Range("F40").Value = 1759
$ws.Range("G40").Value = 65.8
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=85239"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202405/mzD4rhY21715109458100.jpeg"
$ws.Range("F43").Value = 883
$ws.Range("F44").Value = 1042
$ws.Range("F46").Value = 46
$ws.Range("F48").Value = 59
$ws.Range("F49").Value = 68

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 259
$ws.Range("F11").Value = 378
$ws.Range("F30").Value = 297
$ws.Range("F35").Value = 38

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 392
$ws.Range("F7").Value = 1444
$ws.Range("F8").Value = 759
$ws.Range("F9").Value = 376
$ws.Range("F10").Value = 2734
$ws.Range("F11").Value = 262
$ws.Range("F12").Value = 480
$ws.Range("F13").Value = 365
$ws.Range("F14").Value = 1114

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1444
$ws.Range("F3").Value = 759
$ws.Range("F5").Value = 376
$ws.Range("F6").Value = 2734
$ws.Range("F7").Value = 1585
$ws.Range("F8").Value = 769
$ws.Range("F9").Value = 685
$ws.Range("F10").Value = 1281
$ws.Range("F11").Value = 2542
$ws.Range("F12").Value = 1329
$ws.Range("F14").Value = 296
$ws.Range("F15").Value = 2297
$ws.Range("F16").Value = 2000
$ws.Range("F18").Value = 6203
$ws.Range("F19").Value = 111
$ws.Range("F20").Value = 480
$ws.Range("F21").Value = 1205
$ws.Range("F22").Value = 1418
$ws.Range("F23").Value = 1314
$ws.Range("F24").Value = 1182
$ws.Range("F25").Value = 2088
$ws.Range("F29").Value = 679
$ws.Range("F30").Value = 215
$ws.Range("F31").Value = 5223
$ws.Range("F33").Value = 1238
$ws.Range("F34").Value = 3657
$ws.Range("F35").Value = 297
$ws.Range("F36").Value = 1659
$ws.Range("F38").Value = 143
$ws.Range("F39").Value = 956
$ws.Range("F40").Value = 372
$ws.Range("F41").Value = 1759
$ws.Range("F43").Value = 38
$ws.Range("F45").Value = 883
$ws.Range("F46").Value = 1042
$ws.Range("F50").Value = 59
$ws.Range("F51").Value = 68
